# Apply the "lots of discharge data" edit to the stn3 ("sheet2") worksheet.
$wb = $excel.ActiveWorkbook

# stn1, stn3, stn4 -> Worksheets.Item(1), Item(2), Item(3)
$ws3 = $wb.Worksheets.Item(2)   # "stn3" sheet (xl/worksheets/sheet2.xml)

# --- Make "new velocity" label (A14) bold, matching the new "new depth" label style ---
$ws3.Range("A14").Font.Bold = $true

# --- Add the "new depth" section header (row 27), bold, new shared string "new depth" ---
$ws3.Range("A27").Value = "new depth"
$ws3.Range("A27").Font.Bold = $true

# --- Add header row 28 (same labels as rows 2 / 15: x, v, d, segment, Q, Qtotal) ---
$ws3.Range("A28").Value2 = $ws3.Range("A15").Value2
$ws3.Range("B28").Value2 = $ws3.Range("B15").Value2
$ws3.Range("C28").Value2 = $ws3.Range("C15").Value2
$ws3.Range("D28").Value2 = $ws3.Range("D15").Value2
$ws3.Range("E28").Value2 = $ws3.Range("E15").Value2
$ws3.Range("F28").Value2 = $ws3.Range("F15").Value2

# --- Fill the new data block rows 29-37, mirroring rows 16-24 ---
# Column A: same station-depth values as A16:A24
# Column B: literal "new velocity" numbers copied from B16:B24
for ($r = 16; $r -le 24; $r++) {
    $dest = $r + 13
    $ws3.Range("A$dest").Value2 = $ws3.Range("A$r").Value2
    $ws3.Range("B$dest").Value2 = $ws3.Range("B$r").Value2
}

# Column C: depth in cm = corresponding old-table depth (C16:C24) * 2.54
for ($r = 29; $r -le 37; $r++) {
    $src = $r - 13
    $ws3.Range("C$r").Formula = "=C$src*2.54"
}

# Column D: segment midpoints
$ws3.Range("D29").Formula = "=A29"
for ($r = 30; $r -le 37; $r++) {
    $ws3.Range("D$r").Formula = "=(A$r+(A" + ($r + 1) + "-A$r)/2)"
}

# Column E: per-segment discharge
for ($r = 30; $r -le 37; $r++) {
    $prev = $r - 1
    $ws3.Range("E$r").Formula = "=(D$r-D$prev)*(B$r)*C$r"
}

# Column F: total discharge for the new block
$ws3.Range("F29").Formula = "=SUM(E29:E37)"

# --- Update the active window/tab to point at stn3, cell A14 selected (matches diff) ---
$ws3.Activate()
$ws3.Range("A14").Select()
